# Atualizacao rapida de agenda as  9:19:12,18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Rows 2-7: new agenda entries (Giovani / Pendente block) ---
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0643"
$ws.Range("C2").Value = "MegaScan"
$ws.Range("D2").Value = "Sem comunicação de câmeras."
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 21 / Forte: 15"

$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0756"
$ws.Range("C3").Value = "Manoel Correira"
$ws.Range("D3").Value = "Sem comunicação de câmeras, passar pro DDNS pois o cloud tá bem ruim."
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "Pendente"
$ws.Range("H3").ClearContents()

$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0079"
$ws.Range("C4").Value = "Med Center"
$ws.Range("D4").Value = "Sem comunicação de alarmes, internet."
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Pendente"
$ws.Range("H4").ClearContents()

$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0355"
$ws.Range("C5").Value = "Rc Silva"
$ws.Range("D5").Value = "Sem comunicação de alarmes, internet."
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "Pendente"
$ws.Range("H5").ClearContents()
$ws.Rows.Item(5).EntireRow.AutoFit()

$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0098"
$ws.Range("C6").Value = "Localiza"
$ws.Range("D6").Value = "Sem comunicação geral, confirmar se o problema é internet."
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = "Pendente"
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = "Giovani"
$ws.Range("B7").Value = "'0840"
$ws.Range("C7").Value = "Valdemar Amaral"
$ws.Range("D7").Value = "Sem comunicação de alarmes, gprs."
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = "Pendente"
$ws.Range("H7").ClearContents()
$ws.Rows.Item(7).EntireRow.AutoFit()

# --- Rows 8-16: old agenda items removed (cleared), keep formatting ---
$ws.Range("A8:H16").ClearContents()
$ws.Rows.Item(8).EntireRow.AutoFit()
$ws.Rows.Item(9).EntireRow.AutoFit()
$ws.Rows.Item(10).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()
